$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "PDB molecule" column (column D) entirely, shifting
# "PDB filename" and "Is model" left by one column.
$ws.Columns.Item(4).Delete()

# Row 3 (JCVISYN3_0002) referenced the same PDB file as row 2
# (JCVISYN3_0001.pdb via the now-removed "chain" distinction); clear the
# now-duplicate filename value so only the "Is model" flag remains.
$ws.Range("D3").ClearContents()

# Widen the Feature ID / PDB filename columns to fit their contents.
$ws.Columns.Item(3).ColumnWidth = 15.6666666667
$ws.Columns.Item(4).ColumnWidth = 25.6666666667

$ws.Range("F10").Select() | Out-Null
